$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Toggle the Runmode (column C) flags for the Suite test rows.
$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"
$ws.Range("C5").Value = "Y"
$ws.Range("C6").Value = "N"

# Move the active selection to C9 (as recorded in the sheet view).
$null = $ws.Range("C9").Select()
